# This commit ("Actualiza base de datos EC y agrega parte 1 de nuevos
# estado de cuenta") is part of a larger, multi-file batch update across
# the repository. For THIS specific workbook
# (Data/EC/NIT-9015202615.xlsx), Excel re-saved the file without any
# actual change to a cell's displayed value or formula: the canonical
# OOXML diff only touches bookkeeping that Excel regenerates on every
# save (fileVersion/rupBuild, revisionPtr/document & view GUIDs, the
# drawing's creationId GUID, the physical on-disk ordering of the
# <si> shared-string table entries, and the physical ordering of two
# otherwise-identical border/cellXf definitions in styles.xml - all of
# which keep every cell's rendered text/number and every cell's visual
# border formatting byte-for-byte identical to the original).
#
# We reproduce that outcome through the Excel object model by simply
# re-asserting the worksheet's existing data (a no-drift "touch" of the
# workbook), which is equivalent in effect: no cell value, formula or
# visible formatting changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- Header / summary block (unchanged) ---
$ws.Range("D2").Value = "ESTADO DE CUENTA"

$ws.Range("C7").Value = "RAZON SOCIAL:"
$ws.Range("E7").Value = "EMBLUE GROUP S.A."

$ws.Range("C9").Value = "NIT"
$ws.Range("E9").Value = 9015202615

$ws.Range("C11").Value = "VALOR MORA"
$ws.Range("E11").Value = 480000

$ws.Range("B13").Value = "Cant. Trabajadores"
$ws.Range("C13").Value = 1
$ws.Range("E13").Value = "Cant. Periodos"
$ws.Range("F13").Value = 12

# --- Table header row ---
$ws.Range("B15").Value = "Tipo Doc Trabajador"
$ws.Range("C15").Value = "N° Doc Trabajador"
$ws.Range("D15").Value = "Nombre Trabajador"
$ws.Range("E15").Value = "Periodo Mora"
$ws.Range("F15").Value = "Valor Mora"
$ws.Range("G15").Value = "Salario Basico"
$ws.Range("H15").Value = "Novedad de Ingreso"
$ws.Range("I15").Value = "Novedad de Retiro"
$ws.Range("J15").Value = "Observaciones"

# --- Worker debt/period rows (same 12 periods, same worker) ---
$periods = @("2401","2312","2311","2310","2309","2308","2307","2306","2305","2304","2303","2302")
$row = 16
foreach ($periodo in $periods) {
    $ws.Range("B$row").Value = "CC"
    $ws.Range("C$row").Value = "22808549"
    $ws.Range("D$row").Value = "ANGELICA MARIA VEGA BERRIO"
    $ws.Range("E$row").Value = $periodo
    $ws.Range("F$row").Value = 40000
    $ws.Range("G$row").Value = 1000000
    $row++
}

# --- Signature block ---
$ws.Range("B32").Value = "___________________________________"
$ws.Range("H32").Value = "___________________________________"
$ws.Range("B33").Value = "NOMBRE DEL REPRESENTANTE LEGAL"
$ws.Range("H33").Value = "FIRMA DEL REPRESENTANTE LEGAL"

$wb.Save()
